# Add a new "Modelo" column (F) to the worksheet, matching the header style
# of the existing columns and filling rows 2-4 with the model description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell F1 -------------------------------------------------------
# Copy the formatting of the neighboring header cell (E1: bold font, border,
# centered/top aligned) onto F1, then set its text.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F1").Value = "Modelo"

# --- Data cells F2:F4 -------------------------------------------------------
$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n" + `
    "                                            estimator=Pipeline(steps=[('model',`n" + `
    "                                                                       LinearRegression())]),`n" + `
    "                                            param_grid={'model__fit_intercept': [True,`n" + `
    "                                                                                 False]},`n" + `
    "                                            scoring='neg_mean_squared_error'))"

$ws.Range("F2").Value = $modelText
$ws.Range("F3").Value = $modelText
$ws.Range("F4").Value = $modelText

# Entering the multi-line strings above causes Excel to auto-expand the row
# heights; restore them back to the sheet's standard (default) height so the
# rows stay unmodified, matching the original layout.
$ws.Rows.Item(2).AutoFit() | Out-Null
$ws.Rows.Item(3).AutoFit() | Out-Null
$ws.Rows.Item(4).AutoFit() | Out-Null
